# Update odds values on Sheet1 to match the latest FlashScore export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48

# Row 3
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 9

# Row 4
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 7.18

# Row 6
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.85

# Row 7
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 2.3
$ws.Range("L7").Value = 3
$ws.Range("O7").Value = 1.29
$ws.Range("P7").Value = 3.5
$ws.Range("Q7").Value = 1.93
$ws.Range("R7").Value = 1.88
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.75
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 15
$ws.Range("Y7").Value = 11
$ws.Range("AC7").Value = 10
$ws.Range("AD7").Value = 6
$ws.Range("AG7").Value = 8.5
$ws.Range("AH7").Value = 12
$ws.Range("AK7").Value = 19
$ws.Range("AT7").Value = 2.75
$ws.Range("AX7").Value = 13

# Row 12
$ws.Range("G12").Value = 8.25
$ws.Range("H12").Value = 5
$ws.Range("J12").Value = 6.4
$ws.Range("K12").Value = 2.7
$ws.Range("P12").Value = 5.6
$ws.Range("Q12").Value = 1.36
$ws.Range("R12").Value = 2.92
$ws.Range("S12").Value = 1.21
$ws.Range("T12").Value = 3.9
$ws.Range("X12").Value = 75
$ws.Range("Z12").Value = 200
$ws.Range("AA12").Value = 75
$ws.Range("AC12").Value = 10.25
$ws.Range("AD12").Value = 11.25
$ws.Range("AE12").Value = 16
$ws.Range("AG12").Value = 11.25
$ws.Range("AI12").Value = 8.75
$ws.Range("AN12").Value = 9.75
$ws.Range("AP12").Value = 29
$ws.Range("AT12").Value = 3.9
$ws.Range("AU12").Value = 7.2
$ws.Range("AV12").Value = 40
